$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatted style of A2 (bold font + border + center/top alignment)
# onto the new group-name cells A3:A5 so they match the existing header style
# (style index "s=1") without creating new style entries in styles.xml.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3:A5").PasteSpecial(-4122) | Out-Null

# Column A: research-group names (shared strings registered first, matching
# the original workbook author order of entry).
$ws.Range("A2").Value = 'Posgrado en Aprovechamiento de Recursos Hidráulicos'
$ws.Range("A3").Value = 'ESTUDIO Y CONTROL DE LA CONTAMINACIÓN AMBIENTAL - ECCA'
$ws.Range("A4").Value = 'Bioprocesos y Flujos Reactivos'
$ws.Range("A5").Value = 'Limnología y Recursos Hídricos'

# Column B: published-article descriptions.
$ws.Range("B2").Value = '273.- Publicado en revista especializada: Revisión de criterios y metodologías de diseño de redes para el monitoreo de la calidad de agua en ríos  Colombia, Avances En Recursos Hidráulicos ISSN: 0121-5701, 2008 vol:N/A fasc: 18 págs: 67 - 78, DOI:  Autores: A MORENO, FRANCISCO MAURICIO TORO BOTERO, LUIS FERNANDO CARVAJAL SERNA 
 274.- Publicado en revista especializada: Revisión de criterios y metodologías de diseño de redes para el monitoreo de la calidad de agua en ríos  Colombia, Avances En Recursos Hidráulicos ISSN: 0121-5701, 2008 vol:N/A fasc: 18 págs: 67 - 78, DOI:  Autores: LUIS FERNANDO CARVAJAL SERNA, A MORENO, MAURICIO TORO'
$ws.Range("B3").Value = '103.- Publicado en revista especializada: Monitoreo y medición del ajuste del pH del agua tratada del río cauca mediante índices de estabilización  Colombia, Revista U.D.C.A Actualidad and Divulgacion Cientifica ISSN: 0123-4226, 2014 vol:17 fasc: 2 págs: 563 - 575, DOI:  Autores: PATRICIA TORRES LOZADA, KAREN ALEJANDRA BUENO ZABALA, LUIS GERMAN DELGADO 
 176.- Publicado en revista especializada: Propuesta metodológica para localización de estaciones de monitoreo de calidad de agua en redes de distribución utilizando sistemas de información geográfica  Colombia, Revista Facultad de Ingenieria ISSN: 0120-6230, 2009 vol:49 fasc: N/A págs: 129 - 140, DOI:  Autores: CAMILO HERNAN CRUZ VELEZ, CAROLINA MONTOYA PACHONGO, DIANA LOAIZA, PATRICIA TORRES LOZADA, JUAN CARLOS ESCOBAR RIVERA, LUIS GERMAN DELGADO'
$ws.Range("B4").Value = '49.- Publicado en revista especializada: Monitoreo y medición del ajuste del pH del agua tratada del río Cauca mediante índices de estabilización  Colombia, Revista U.D.C.A Actualidad and Divulgacion Cientifica ISSN: 0123-4226, 2014 vol:17 fasc: 2 págs: 563 - 575, DOI:  Autores: KAREN ALEJANDRA BUENO ZABALA, PATRICIA TORRES LOZADA, LUIS GERMAN DELGADO CABRERA'
$ws.Range("B5").Value = '38.- Corto (Resumen): Monitoreo de la asociacion de especies icticas, aguas arriba del embalse Amani, rio La Miel aguas abajode la presa de la central hidroelectrica Miel I- y afluentes principales (Rios Manso y Samana) duante 2008, Colombia  Colombia, Actualidades Biológicas ISSN: 0304-3584, 2009 vol:31 fasc: N/A págs: 54 - 54, DOI:  Autores: MARIA ISABEL RIOS PULGARIN, ALEXANDRA ARANGO ROJAS, ESTEBAN PELAEZ SANCHEZ 
 44.- Publicado en revista especializada: Monitoreo de la asociaciones de especies ícticas, aguas arriba del embalse Amaní, Río La Miel aguas abajo de la presa de la Hidroeléctrica Miel I y afluentes principales (Ríos Manso y Samaná) durante 2008  Colombia, Actualidades Biológicas ISSN: 0304-3584, 2009 vol:31 fasc: N/A págs: 54 - 54, DOI:  Autores: MARIA ISABEL RIOS PULGARIN, GLORIA ALEXANDRA ARANGO ROJAS, ESTEBAN PELAEZ SANCHEZ'
